$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VCO")
$ws.Activate()

# Add the new SI570 row (row 8) to the VCO / CVO comparison table.
# Values are written in D, C, B, A order so that the workbook's shared
# string table receives the new unique strings in the same order as the
# target file (url, supply voltage, frequency range, part name).
$ws.Range("D8").Value = "https://www.silabs.com/documents/public/data-sheets/si570.pdf"

$ws.Range("C8").Value = "1.8, 2.5, or 3.3 V supply"
$ws.Range("C8").Style = "Normal"

$ws.Range("B8").Value = " 10 MHz to 1400 MHz"
$ws.Range("B8").WrapText = $true

$ws.Range("A8").Value = "Si570"

# Move/keep the selection on the newly filled cell, like in the saved file.
$ws.Range("A8").Select()
